# edit.ps1 -- apply "Updated to December 2022" commit
# Adds Sep/Okt/Nov/Dec 2022 monthly data (rows 12-15) to the "Electricity 2022"
# sheet, recomputes the dependent totals/formulas, swaps the two "2021"-labelled
# summary captions for fresh "2022" text, clears the now-stale Y20 helper
# formula, adds the H18 "production vs. estimate" percentage, and restores the
# expected selections on the 2021/2022 sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Electricity 2020"
$ws2 = $wb.Worksheets.Item(2)   # "Electricity 2021"
$ws3 = $wb.Worksheets.Item(3)   # "Electricity 2022"

# --- Fill in the newly-reported months (Sep/Okt/Nov/Dec 2022), rows 12-15 ---

# Row 12
$ws3.Range("C12").Value = 1219556
$ws3.Range("D12").Value = 765749
$ws3.Range("E12").Value = 602849
$ws3.Range("F12").Value = 616707
$ws3.Range("G12").Value = 1382456
$ws3.Range("H12").Formula = "=(G12/1000)/B12"
$ws3.Range("I12").Value = 598
$ws3.Range("J12").Value = 763
$ws3.Range("K12").Formula = "=(G12/1000)+I12-J12"
$ws3.Range("L12").Formula = "=K12-I12"
$ws3.Range("M12").Formula = "=L12/K12"
$ws3.Range("N12").Value = 0
$ws3.Range("O12").Value = 452.5
$ws3.Range("P12").Value = 82.29
$ws3.Range("Q12").Value = 17
$ws3.Range("R12").Value = 45
$ws3.Range("S12").Value = 346.05
$ws3.Range("T12").Value = 6.9
$ws3.Range("U12").Value = 306.43
$ws3.Range("V12").Formula = "=P12*I12/100"
$ws3.Range("W12").Formula = "=U12*I12/100"
$ws3.Range("X12").Formula = "=(K12*(U12+Q12+R12)/100)+N12+O12"
$ws3.Range("Y12").Formula = "=(O12+((Q12+R12)*I12/100)+V12)-Z12"
$ws3.Range("Z12").Formula = "=J12*(S12+T12)/100"
$ws3.Range("AA12").Formula = "=X12-Y12"
$ws3.Range("AB12").Formula = "=J12*0.6"
$ws3.Range("AC12").Formula = "=AA12+AB12"
$ws3.Range("AD12").Formula = "=AC12/(G12/1000)"
$ws3.Range("AE12").Formula = "=(P12+Q12+R12)/100"
$ws3.Range("AF12").Value = 300

# Row 13
$ws3.Range("C13").Value = 1304831
$ws3.Range("D13").Value = 273151
$ws3.Range("E13").Value = 789423
$ws3.Range("F13").Value = 515408
$ws3.Range("G13").Value = 788559
$ws3.Range("H13").Formula = "=(G13/1000)/B13"
$ws3.Range("I13").Value = 784
$ws3.Range("J13").Value = 270
$ws3.Range("K13").Formula = "=(G13/1000)+I13-J13"
$ws3.Range("L13").Formula = "=K13-I13"
$ws3.Range("M13").Formula = "=L13/K13"
$ws3.Range("N13").Value = 0
$ws3.Range("O13").Value = 452.5
$ws3.Range("P13").Value = 45.35
$ws3.Range("Q13").Value = 17
$ws3.Range("R13").Value = 45
$ws3.Range("S13").Value = 146.47
$ws3.Range("T13").Value = 6.9
$ws3.Range("U13").Value = 114.91
$ws3.Range("V13").Formula = "=P13*I13/100"
$ws3.Range("W13").Formula = "=U13*I13/100"
$ws3.Range("X13").Formula = "=(K13*(U13+Q13+R13)/100)+N13+O13"
$ws3.Range("Y13").Formula = "=(O13+((Q13+R13)*I13/100)+V13)-Z13"
$ws3.Range("Z13").Formula = "=J13*(S13+T13)/100"
$ws3.Range("AA13").Formula = "=X13-Y13"
$ws3.Range("AB13").Formula = "=J13*0.6"
$ws3.Range("AC13").Formula = "=AA13+AB13"
$ws3.Range("AD13").Formula = "=AC13/(G13/1000)"
$ws3.Range("AE13").Formula = "=(P13+Q13+R13)/100"
$ws3.Range("AF13").Value = 300

# Row 14
$ws3.Range("C14").Value = 1410228
$ws3.Range("D14").Value = 29350
$ws3.Range("E14").Value = 1281077
$ws3.Range("F14").Value = 131674
$ws3.Range("G14").Value = 158501
$ws3.Range("H14").Formula = "=(G14/1000)/B14"
$ws3.Range("I14").Value = 1276
$ws3.Range("J14").Value = 26
$ws3.Range("K14").Formula = "=(G14/1000)+I14-J14"
$ws3.Range("L14").Formula = "=K14-I14"
$ws3.Range("M14").Formula = "=L14/K14"
$ws3.Range("N14").Value = 0
$ws3.Range("O14").Value = 452.5
$ws3.Range("P14").Value = 124.53
$ws3.Range("Q14").Value = 17
$ws3.Range("R14").Value = 45
$ws3.Range("S14").Value = 102.21
$ws3.Range("T14").Value = 6.9
$ws3.Range("U14").Value = 197.51
$ws3.Range("V14").Formula = "=P14*I14/100"
$ws3.Range("W14").Formula = "=U14*I14/100"
$ws3.Range("X14").Formula = "=(K14*(U14+Q14+R14)/100)+N14+O14"
$ws3.Range("Y14").Formula = "=(O14+((Q14+R14)*I14/100)+V14)-Z14"
$ws3.Range("Z14").Formula = "=J14*(S14+T14)/100"
$ws3.Range("AA14").Formula = "=X14-Y14"
$ws3.Range("AB14").Formula = "=J14*0.6"
$ws3.Range("AC14").Formula = "=AA14+AB14"
$ws3.Range("AD14").Formula = "=AC14/(G14/1000)"
$ws3.Range("AE14").Formula = "=(P14+Q14+R14)/100"
$ws3.Range("AF14").Value = 300

# Row 15
$ws3.Range("C15").Value = 1687518
$ws3.Range("D15").Value = 13500
$ws3.Range("E15").Value = 1668805
$ws3.Range("F15").Value = 23707
$ws3.Range("G15").Value = 32213
$ws3.Range("H15").Formula = "=(G15/1000)/B15"
$ws3.Range("I15").Value = 1666
$ws3.Range("J15").Value = 5.4
$ws3.Range("K15").Formula = "=(G15/1000)+I15-J15"
$ws3.Range("L15").Formula = "=K15-I15"
$ws3.Range("M15").Formula = "=L15/K15"
$ws3.Range("N15").Value = 0
$ws3.Range("O15").Value = 452.5
$ws3.Range("P15").Value = 258.41
$ws3.Range("Q15").Value = 17
$ws3.Range("R15").Value = 45
$ws3.Range("S15").Value = 235.83
$ws3.Range("T15").Value = 6.9
$ws3.Range("U15").Value = 371.36
$ws3.Range("V15").Formula = "=P15*I15/100"
$ws3.Range("W15").Formula = "=U15*I15/100"
$ws3.Range("X15").Formula = "=(K15*(U15+Q15+R15)/100)+N15+O15"
$ws3.Range("Y15").Formula = "=(O15+((Q15+R15)*I15/100)+V15)-Z15"
$ws3.Range("Z15").Formula = "=J15*(S15+T15)/100"
$ws3.Range("AA15").Formula = "=X15-Y15"
$ws3.Range("AB15").Formula = "=J15*0.6"
$ws3.Range("AC15").Formula = "=AA15+AB15"
$ws3.Range("AD15").Formula = "=AC15/(G15/1000)"
$ws3.Range("AE15").Formula = "=(P15+Q15+R15)/100"
$ws3.Range("AF15").Value = 300


# --- H18: new "Actual vs. estimated production" percentage (0.0%) ---
$ws3.Range("H18").Formula = "=G18/B18"
$ws3.Range("H18").NumberFormat = "0.0%"

# --- Y20 no longer shows the stale "Y18-AB18" helper figure ---
$ws3.Range("Y20").ClearContents() | Out-Null

# --- Refresh the two summary captions from "...2021" to "...2022" ---
$ws3.Range("C24").Value = "Total value of solar cell production 2022"
$ws3.Range("C27").Value = "Total value of solar cells after intrest 2022"

# --- Restore expected selections / scroll position on sheets 2021 and 2022 ---
$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws2.Range("Y20").Select() | Out-Null

$ws3.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws3.Range("C28").Select() | Out-Null

